# middleman-bot/middleman-bot/lifts.xlsx
# user_data & bot_data keys from string to enums + Ready state +
# Add per user signaling and linking : progress

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook-level: rename "temp" -> "follows"
# ---------------------------------------------------------------------
$wsLifts   = $wb.Worksheets.Item(1)
$wsFollows = $wb.Worksheets.Item(2)
$wsFollows.Name = "follows"

# Best-effort restore of the saved window geometry (bookViews/workbookView).
$win = $wb.Windows.Item(1)
$win.Height = 28590
$win.Width  = 16440
$win.Left   = 28680
$win.Top    = -8835

# =======================================================================
# Sheet "lifts"
# =======================================================================

# -- Column widths -------------------------------------------------------
# B gets a dedicated wider column, C:G keep the default width but become
# their own explicit block (matches the new <cols> split in the diff).
$wsLifts.Columns.Item(2).ColumnWidth = 15.5
$wsLifts.Range("C1:G1").EntireColumn.ColumnWidth = 8.3

# -- New "USERS" header column (F1) --------------------------------------
$wsLifts.Range("F1").Value = "USERS"
$wsLifts.Range("F1").Font.Bold = $true
$wsLifts.Range("F1").HorizontalAlignment = -4108  # xlCenter
$wsLifts.Range("F1").VerticalAlignment   = -4108  # xlCenter

# -- Rows 2-4 lose their centred "header-like" styling -------------------
$wsLifts.Range("A2:E4").Style = "Normal"

# -- Cell value changes ----------------------------------------------------
$wsLifts.Range("E3").Value = "N"
$wsLifts.Range("E4").Value = "Nn"

$wsLifts.Range("B5").Value = "NONE"
$wsLifts.Range("E5").Value = "N"

$wsLifts.Range("B6").Value = "OPENING"
$wsLifts.Range("E6").Value = "Note"

$wsLifts.Range("B7").Value = "READY"

# =======================================================================
# Sheet "follows" (formerly "temp")
# =======================================================================
$wsFollows.Range("A1").Value = "S"
$wsFollows.Range("B1").Value = "046G"
$wsFollows.Range("C1").Value = "047G"

$wsFollows.Range("A1:C1").Font.Bold = $true
$wsFollows.Range("A1:C1").HorizontalAlignment = -4108  # xlCenter
$wsFollows.Range("A1:C1").VerticalAlignment   = -4108  # xlCenter

$wsFollows.PageSetup.PaperSize   = 9  # xlPaperA4
$wsFollows.PageSetup.Orientation = 1  # xlPortrait

# Widen the explicit <cols> block out to column G (matches the new
# dimension/<cols> split in the diff).
$wsFollows.Range("E1:G1").EntireColumn.ColumnWidth = 8.3

$wsFollows.Activate()
$wsFollows.Range("D6").Select() | Out-Null

# -- Selection (lifts is the tab that stays active/selected) --------------
$wsLifts.Activate()
$wsLifts.Range("G12").Select() | Out-Null

Write-Output "edit complete"
